$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the header date (A1): 45406 -> 45436
$ws.Range("A1").Value = 45436

# Update the price column (D33:D37)
$ws.Range("D33").Value = 767.647
$ws.Range("D34").Value = 1139.001
$ws.Range("D35").Value = 1427.198
$ws.Range("D36").Value = 1718.204
$ws.Range("D37").Value = 1878.96
